$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

# Update the "Promedio" text in E2 to include the LaTeX \href certificate link,
# and remove the separate "Certificado: ..." row that used to live on row 3.
$ws.Range("E2").Value = "Promedio: 97/100 (ver \href{https://www.coursera.org/account/accomplishments/verify/DC7ULMJ3CZWM}{certificado})"

# Delete the row that only contained the old "Certificado: ..." text (row 3),
# shifting the Dundee entry (previously row 4) up into row 3.
$ws.Range("A3:XFD3").EntireRow.Delete()

# Match the saved selection/view state from the edit (row 3 fully selected,
# scrolled so column C is the left-most visible column).
$ws.Range("A3:XFD3").Select()
$excel.ActiveWindow.ScrollColumn = 3
